# Rename the inventory header row to the lowercase/short titles used
# across the comparator tool, and refresh the formatting/selection that
# Excel re-stamped on A1 / D7 while the sheet was being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header titles (row 1) - same cells/columns, updated wording.
$ws.Range("A1").Value = "part number"
$ws.Range("B1").Value = "quantity"
$ws.Range("C1").Value = "UM"
$ws.Range("D1").Value = "value mxn"

# A1 and D7 pick up an explicit font stamp (re-applied formatting) as
# part of the edit; touch their font so the style gets recorded.
$ws.Range("A1").Font.ThemeColor = 1
$ws.Range("D7").Font.ThemeColor = 1

# Leave the cursor on D7, matching the end of the editing session.
[void]$ws.Range("D7").Select()
